# documentation/variables.xlsx — "math checks: locate differences between
# 2016 script and manual computation"
#
# Row 6 (the "OA" / older-adults universe row) was pointing at the wrong
# Census table for its Count/Universe columns (it had been copy-pasted from
# the LEP/"language spoken at home" row). Point B6/C6 back at the
# S0101_C01_001 (AGE AND SEX, total population) variable that the other
# "Count" column entries use, and correct the Percent-column description
# (E6) to the DP05_0029E "65 years and over" variable actually referenced
# by the rest of that row (D6/F6/G6 already pointed at the 65-and-over
# numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "S0101_C01_001"
$ws.Range("C6").Value = "S0101_C01_001" + [char]10 + "AGE AND SEX" + [char]10 + "Estimate!!Total!!Total population"

$ws.Range("E6").Value = "DP05_0029E" + [char]10 + "ACS DEMOGRAPHIC AND HOUSING ESTIMATES" + [char]10 + " Estimate!!SEX AND AGE!!Total population!!65 years and over"

# The previous selection/scroll position (topLeftCell A7, active cell D4)
# was stale; reselect G8 as the last-touched cell with the sheet scrolled
# back to the top.
$ws.Range("G8").Select()
